$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.470.10"
$ws.Range("E2").Value = "  -0.17%  "

$ws.Range("D3").Value = "1.841.10"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'261.88"
$ws.Range("E5").Value = "  -0.48%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").Value = "'0.5301"
$ws.Range("E7").Value = "  +1.40%  "

$ws.Range("D8").Value = "'0.3065"
$ws.Range("E8").Value = "  -5.32%  "

$ws.Range("D9").Value = "'0.06892"

$ws.Range("D10").Value = "'18.38"
$ws.Range("E10").Value = "  -1.76%  "

$ws.Range("D11").Value = "'0.07803"
$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("D12").Value = "'0.7548"
$ws.Range("E12").Value = "  -2.08%  "

$ws.Range("D13").Value = "1.840.54"
$ws.Range("E13").Value = "  -1.21%  "

$ws.Range("D14").Value = "'89.47"
$ws.Range("E14").Value = "  +1.29%  "

$ws.Range("D15").Value = "'5.015"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").Value = "'13.99"
$ws.Range("E17").Value = "  +0.52%  "

$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007930"
$ws.Range("E19").Value = "  +0.11%  "

$ws.Range("D20").Value = "26.492.99"
$ws.Range("E20").Value = "  -0.26%  "

$ws.Range("D21").Value = "'4.614"
$ws.Range("E21").Value = "  -0.07%  "

$ws.Range("D22").Value = "'5.979"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").Value = "'9.298"
$ws.Range("E23").Value = "  -1.26%  "

$ws.Range("D24").Value = "'142.20"
$ws.Range("E24").Value = "  -0.37%  "

$ws.Range("D25").Value = "'2.188"
$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("D26").Value = "'1.693"
$ws.Range("E26").Value = "  +0.90%  "

$ws.Range("D27").Value = "'17.00"
$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D28").Value = "'111.28"
$ws.Range("E28").Value = "  -0.29%  "

$ws.Range("D29").Value = "'4.264"
$ws.Range("E29").Value = "  +2.34%  "

$ws.Range("D30").Value = "'0.08784"
$ws.Range("E30").Value = "  +0.46%  "

$ws.Range("D31").Value = "'4.078"
$ws.Range("E31").Value = "  -0.78%  "

$ws.Range("D32").Value = "'0.04814"
$ws.Range("E32").Value = "  -0.13%  "

$ws.Range("D33").Value = "'2.930"
$ws.Range("E33").Value = "  +1.90%  "

$ws.Range("D34").Value = "'0.7288"
$ws.Range("E34").Value = "  +1.91%  "

$ws.Range("D35").Value = "'1.131"
$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("D36").Value = "'3.102"
$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").Value = "'2.306"
$ws.Range("E37").Value = "  +5.44%  "

$ws.Range("D38").Value = "'0.01716"
$ws.Range("E38").Value = "  -3.82%  "

$ws.Range("D39").Value = "'0.4788"
$ws.Range("E39").Value = "  -0.97%  "

$ws.Range("D40").Value = "'0.9020"
$ws.Range("E40").Value = "  +0.72%  "

$ws.Range("D41").Value = "'107.89"
$ws.Range("E41").Value = "  -4.02%  "

$ws.Range("D42").Value = "'5.877"
$ws.Range("E42").Value = "  -2.75%  "

$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").Value = "'7.450"
$ws.Range("E44").Value = "  -2.34%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'9.096"
$ws.Range("E45").Value = "  +0.73%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.4126"
$ws.Range("E46").Value = "  -0.92%  "

$ws.Range("D47").Value = "'0.1237"
$ws.Range("E47").Value = "  +0.79%  "

$ws.Range("D48").Value = "'34.85"
$ws.Range("E48").Value = "  -0.21%  "

$ws.Range("D49").Value = "'0.8980"
$ws.Range("E49").Value = "  +1.58%  "

$ws.Range("D50").Value = "'0.05797"

$ws.Range("D51").Value = "'60.11"
$ws.Range("E51").Value = "  +0.26%  "
